$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 574749
$ws.Range("E2").Value = 57876
$ws.Range("F2").Value = 57876
$ws.Range("G2").Value = 42293
$ws.Range("H2").Value = 27990
$ws.Range("I2").Value = 26869
$ws.Range("J2").Value = 1121
$ws.Range("K2").Value = 1637083
$ws.Range("L2").Value = 1088833
$ws.Range("M2").Value = 548250
$ws.Range("N2").Value = 536013
$ws.Range("O2").Value = 12237
$ws.Range("P2").Value = 32098
$ws.Range("Q2").Value = 120457
$ws.Range("R2").Value = -144603
$ws.Range("S2").Value = 19852
$ws.Range("T2").Value = 145475
$ws.Range("U2").Value = -25018
$ws.Range("V2").Value = 635938
$ws.Range("W2").Value = 10.07
$ws.Range("X2").Value = 4.87
$ws.Range("Y2").Value = 5.17
$ws.Range("Z2").Value = 1.75
$ws.Range("AA2").Value = 198.6
$ws.Range("AB2").Value = 1162.02
$ws.Range("AC2").Value = 4185
$ws.Range("AD2").Value = 10.2
$ws.Range("AE2").Value = 83496
$ws.Range("AF2").Value = 0.51
$ws.Range("AG2").Value = 500
$ws.Range("AH2").Value = 1.17
$ws.Range("AI2").Value = 11.95
$ws.Range("AJ2").Value = 641964077

# Row 3
$ws.Range("D3").Value = 589577
$ws.Range("E3").Value = 113467
$ws.Range("F3").Value = 113467
$ws.Range("G3").Value = 186558
$ws.Range("H3").Value = 134164
$ws.Range("I3").Value = 132891
$ws.Range("J3").Value = 1272
$ws.Range("K3").Value = 1752574
$ws.Range("L3").Value = 1073149
$ws.Range("M3").Value = 679425
$ws.Range("N3").Value = 666345
$ws.Range("O3").Value = 13080
$ws.Range("P3").Value = 32098
$ws.Range("Q3").Value = 169431
$ws.Range("R3").Value = -97740
$ws.Range("S3").Value = -52066
$ws.Range("T3").Value = 140499
$ws.Range("U3").Value = 28932
$ws.Range("V3").Value = 594129
$ws.Range("W3").Value = 19.25
$ws.Range("X3").Value = 22.76
$ws.Range("Y3").Value = 22.11
$ws.Range("Z3").Value = 7.92
$ws.Range("AA3").Value = 157.95
$ws.Range("AB3").Value = 1564.84
$ws.Range("AC3").Value = 20701
$ws.Range("AD3").Value = 2.42
$ws.Range("AE3").Value = 103798
$ws.Range("AF3").Value = 0.48
$ws.Range("AG3").Value = 3100
$ws.Range("AH3").Value = 6.2
$ws.Range("AI3").Value = 14.98
$ws.Range("AJ3").Value = 641964077

# Row 4
$ws.Range("D4").Value = 601904
$ws.Range("E4").Value = 120016
$ws.Range("F4").Value = 120016
$ws.Range("G4").Value = 105135
$ws.Range("H4").Value = 71483
$ws.Range("I4").Value = 70486
$ws.Range("J4").Value = 997
$ws.Range("K4").Value = 1778370
$ws.Range("L4").Value = 1047865
$ws.Range("M4").Value = 730505
$ws.Range("N4").Value = 717237
$ws.Range("O4").Value = 13269
$ws.Range("P4").Value = 32098
$ws.Range("Q4").Value = 165206
$ws.Range("R4").Value = -96459
$ws.Range("S4").Value = -76375
$ws.Range("T4").Value = 120288
$ws.Range("U4").Value = 44918
$ws.Range("V4").Value = 541804
$ws.Range("W4").Value = 19.94
$ws.Range("X4").Value = 11.88
$ws.Range("Y4").Value = 10.19
$ws.Range("Z4").Value = 4.05
$ws.Range("AA4").Value = 143.44
$ws.Range("AB4").Value = 1721.37
$ws.Range("AC4").Value = 10980
$ws.Range("AD4").Value = 4.01
$ws.Range("AE4").Value = 111725
$ws.Range("AF4").Value = 0.39
$ws.Range("AG4").Value = 1980
$ws.Range("AH4").Value = 4.49
$ws.Range("AI4").Value = 18.03
$ws.Range("AJ4").Value = 641964077

# Row 5
$ws.Range("D5").Value = 598149
$ws.Range("E5").Value = 49532
$ws.Range("F5").Value = 49532
$ws.Range("G5").Value = 36142
$ws.Range("H5").Value = 14414
$ws.Range("I5").Value = 12987
$ws.Range("J5").Value = 1427
$ws.Range("K5").Value = 1817889
$ws.Range("L5").Value = 1088243
$ws.Range("M5").Value = 729646
$ws.Range("N5").Value = 716814
$ws.Range("O5").Value = 12832
$ws.Range("P5").Value = 32098
$ws.Range("Q5").Value = 112499
$ws.Range("R5").Value = -126067
$ws.Range("S5").Value = 7456
$ws.Range("T5").Value = 125360
$ws.Range("U5").Value = -12861
$ws.Range("V5").Value = 551657
$ws.Range("W5").Value = 8.279999999999999
$ws.Range("X5").Value = 2.41
$ws.Range("Y5").Value = 1.81
$ws.Range("Z5").Value = 0.8
$ws.Range("AA5").Value = 149.15
$ws.Range("AB5").Value = 1727.45
$ws.Range("AC5").Value = 2023
$ws.Range("AD5").Value = 18.86
$ws.Range("AE5").Value = 111660
$ws.Range("AF5").Value = 0.34
$ws.Range("AG5").Value = 790
$ws.Range("AH5").Value = 2.07
$ws.Range("AI5").Value = 39.05
$ws.Range("AJ5").Value = 641964077

# Row 6
$ws.Range("D6").Value = 606276
$ws.Range("E6").Value = -2080
$ws.Range("F6").Value = -2080
$ws.Range("G6").Value = -20008
$ws.Range("H6").Value = -11745
$ws.Range("I6").Value = -13146
$ws.Range("K6").Value = 1852491
$ws.Range("L6").Value = 1141563
$ws.Range("M6").Value = 710928
$ws.Range("N6").Value = 697439
$ws.Range("P6").Value = 32098
$ws.Range("Q6").Value = 66801
$ws.Range("R6").Value = -130143
$ws.Range("S6").Value = 53017
$ws.Range("T6").Value = 122669
$ws.Range("U6").Value = -55867
$ws.Range("V6").Value = 613184
$ws.Range("W6").Value = -0.34
$ws.Range("X6").Value = -1.94
$ws.Range("Y6").Value = -1.86
$ws.Range("Z6").Value = -0.64
$ws.Range("AA6").Value = 160.57
$ws.Range("AB6").Value = 1669.8
$ws.Range("AC6").Value = -2048
$ws.Range("AD6").Value = -16.16
$ws.Range("AE6").Value = 108641
$ws.Range("AF6").Value = 0.3
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 641964077

# Row 7
$ws.Range("D7").Value = 591447
$ws.Range("E7").Value = -3396
$ws.Range("G7").Value = -21948
$ws.Range("H7").Value = -16079
$ws.Range("I7").Value = -16875
$ws.Range("K7").Value = 1911447
$ws.Range("L7").Value = 1217756
$ws.Range("M7").Value = 694068
$ws.Range("N7").Value = 679959
$ws.Range("P7").Value = 32099
$ws.Range("Q7").Value = 96439
$ws.Range("R7").Value = -139643
$ws.Range("S7").Value = 42376
$ws.Range("T7").Value = 130918
$ws.Range("U7").Value = -41978
$ws.Range("W7").Value = -0.57
$ws.Range("X7").Value = -2.72
$ws.Range("Y7").Value = -2.45
$ws.Range("Z7").Value = -0.85
$ws.Range("AA7").Value = 175.45
$ws.Range("AC7").Value = -2629
$ws.Range("AD7").Value = -9.640000000000001
$ws.Range("AE7").Value = 105919
$ws.Range("AF7").Value = 0.24
$ws.Range("AG7").Value = 0
$ws.Range("AH7").Value = 0

# Row 8
$ws.Range("D8").Value = 600722
$ws.Range("E8").Value = 30557
$ws.Range("G8").Value = 17357
$ws.Range("H8").Value = 12486
$ws.Range("I8").Value = 11948
$ws.Range("K8").Value = 1954710
$ws.Range("L8").Value = 1248916
$ws.Range("M8").Value = 705398
$ws.Range("N8").Value = 690470
$ws.Range("P8").Value = 32099
$ws.Range("Q8").Value = 125076
$ws.Range("R8").Value = -144547
$ws.Range("S8").Value = 19271
$ws.Range("T8").Value = 142456
$ws.Range("U8").Value = -15490
$ws.Range("W8").Value = 5.09
$ws.Range("X8").Value = 2.08
$ws.Range("Y8").Value = 1.74
$ws.Range("Z8").Value = 0.65
$ws.Range("AA8").Value = 177.05
$ws.Range("AC8").Value = 1861
$ws.Range("AD8").Value = 13.62
$ws.Range("AE8").Value = 107556
$ws.Range("AF8").Value = 0.24
$ws.Range("AG8").Value = 428
$ws.Range("AH8").Value = 1.69
$ws.Range("AI8").Value = 22.98

# Row 9
$ws.Range("D9").Value = 611858
$ws.Range("E9").Value = 35607
$ws.Range("G9").Value = 23505
$ws.Range("H9").Value = 16840
$ws.Range("I9").Value = 16394
$ws.Range("K9").Value = 1991744
$ws.Range("L9").Value = 1271790
$ws.Range("M9").Value = 718955
$ws.Range("N9").Value = 703558
$ws.Range("P9").Value = 32099
$ws.Range("Q9").Value = 130809
$ws.Range("R9").Value = -146138
$ws.Range("S9").Value = 14722
$ws.Range("T9").Value = 143310
$ws.Range("U9").Value = -14059
$ws.Range("W9").Value = 5.82
$ws.Range("X9").Value = 2.75
$ws.Range("Y9").Value = 2.35
$ws.Range("Z9").Value = 0.85
$ws.Range("AA9").Value = 176.89
$ws.Range("AC9").Value = 2554
$ws.Range("AD9").Value = 9.93
$ws.Range("AE9").Value = 109595
$ws.Range("AF9").Value = 0.23
$ws.Range("AG9").Value = 502
$ws.Range("AH9").Value = 1.98
$ws.Range("AI9").Value = 19.67

# Remove cells that no longer exist after the edit
$ws.Range("AG6").ClearContents()
$ws.Range("AH6").ClearContents()
$ws.Range("AI7").ClearContents()
